$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change
$ws.Range("G1").Value = "table_header_position"

# Rows 2-13: "center right position" header descriptor -> new "top right position" text
$centerRightText = "`"Part Number - Can be found on the top right position of the page`""
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 7).Value = $centerRightText
}

# Rows 14-25: "right side" descriptor text update
$rightSideText = "`"Part Number Description Dimensions Power Lumens Colour Temp. - Can be found on the right side of the page`""
for ($r = 14; $r -le 25; $r++) {
    $ws.Cells.Item($r, 7).Value = $rightSideText
}

# Power/Lumens value swaps
# Row 3: D3/E3 empty -> "16W"/"1600lm"
$ws.Range("D3").Value = "16W"
$ws.Range("E3").Value = "1600lm"

# Row 11: D11/E11 empty -> "16W"/"1600lm"
$ws.Range("D11").Value = "16W"
$ws.Range("E11").Value = "1600lm"

# Row 12: D12/E12 "16W"/"1600lm" -> empty
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""

# Row 13: D13/E13 empty -> "16W"/"1600lm"
$ws.Range("D13").Value = "16W"
$ws.Range("E13").Value = "1600lm"
